# Updates cryptos list data: prices and 1h volume percentages,
# plus a rank swap between RenderToken and Bittensor (rows 37-38).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'95.949.27"
$ws.Range("E2").Value = "'  +0.70%  "
$ws.Range("D3").Value = "'3.564.86"
$ws.Range("E3").Value = "'  -1.62%  "
$ws.Range("D5").Value = "'239.32"
$ws.Range("E5").Value = "'  +0.56%  "
$ws.Range("D6").Value = "'653.49"
$ws.Range("E6").Value = "'  +0.09%  "
$ws.Range("E7").Value = "'  +9.94%  "
$ws.Range("D8").Value = "'0.406"
$ws.Range("E8").Value = "'  +0.36%  "
$ws.Range("E9").Value = "'  -0.01%  "
$ws.Range("D10").Value = "'1.05"
$ws.Range("E10").Value = "'  +5.56%  "
$ws.Range("D11").Value = "'3.563.05"
$ws.Range("E11").Value = "'  -1.62%  "
$ws.Range("D12").Value = "'43.09"
$ws.Range("E12").Value = "'  +1.16%  "
$ws.Range("E13").Value = "'  +0.44%  "
$ws.Range("E14").Value = "'  +0.67%  "
$ws.Range("D15").Value = "'4.227.92"
$ws.Range("E15").Value = "'  -2.09%  "
$ws.Range("D16").Value = "'95.855.91"
$ws.Range("E16").Value = "'  +0.67%  "
$ws.Range("E17").Value = "'  +1.75%  "
$ws.Range("D18").Value = "'3.556.80"
$ws.Range("E18").Value = "'  -1.80%  "
$ws.Range("D19").Value = "'7.75"
$ws.Range("E19").Value = "'  -2.18%  "
$ws.Range("D20").Value = "'12.65"
$ws.Range("E20").Value = "'  -2.66%  "
$ws.Range("D21").Value = "'17.70"
$ws.Range("E21").Value = "'  -1.53%  "
$ws.Range("D22").Value = "'0.510"
$ws.Range("E22").Value = "'  +6.14%  "
$ws.Range("D23").Value = "'3.39"
$ws.Range("E23").Value = "'  -5.25%  "
$ws.Range("D24").Value = "'500.64"
$ws.Range("E24").Value = "'  -1.24%  "
$ws.Range("D25").Value = "'6.91"
$ws.Range("E25").Value = "'  +4.53%  "
$ws.Range("E26").Value = "'  +0.34%  "
$ws.Range("D27").Value = "'95.67"
$ws.Range("E27").Value = "'  -0.36%  "
$ws.Range("D28").Value = "'12.81"
$ws.Range("E28").Value = "'  +0.77%  "
$ws.Range("D29").Value = "'3.755.40"
$ws.Range("E29").Value = "'  -1.55%  "
$ws.Range("E30").Value = "'  +9.47%  "
$ws.Range("D31").Value = "'2.99"
$ws.Range("E31").Value = "'  -4.01%  "
$ws.Range("D32").Value = "'11.33"
$ws.Range("E32").Value = "'  +0.18%  "
$ws.Range("E33").Value = "'  +0.13%  "
$ws.Range("E34").Value = "'  +2.92%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "'  +0.77%  "
$ws.Range("D36").Value = "'31.20"
$ws.Range("E36").Value = "'  -2.93%  "
$ws.Range("B37").Value = "'Bittensor"
$ws.Range("C37").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "'616.25"
$ws.Range("E37").Value = "'  +7.76%  "
$ws.Range("B38").Value = "'RenderToken"
$ws.Range("C38").Value = "'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").Value = "'8.71"
$ws.Range("E38").Value = "'  +6.88%  "
$ws.Range("D39").Value = "'0.561"
$ws.Range("E39").Value = "'  +0.38%  "
$ws.Range("D40").Value = "'1.64"
$ws.Range("E40").Value = "'  +10.77%  "
$ws.Range("E42").Value = "'  +0.07%  "
$ws.Range("D43").Value = "'0.899"
$ws.Range("E43").Value = "'  -3.28%  "
$ws.Range("E44").Value = "'  +5.31%  "
$ws.Range("D45").Value = "'5.69"
$ws.Range("E45").Value = "'  +0.33%  "
$ws.Range("D46").Value = "'23.52"
$ws.Range("E46").Value = "'  -0.74%  "
$ws.Range("D47").Value = "'0.0420"
$ws.Range("E47").Value = "'  +1.90%  "
$ws.Range("D48").Value = "'2.26"
$ws.Range("E48").Value = "'  +1.47%  "
$ws.Range("D49").Value = "'33.53"
$ws.Range("E49").Value = "'  -4.92%  "
$ws.Range("D50").Value = "'3.50"
$ws.Range("E50").Value = "'  -0.44%  "
$ws.Range("D51").Value = "'8.16"
$ws.Range("E51").Value = "'  +1.44%  "
